# Append four new "befor -> after" rows to the normalization table on
# sheet "normal" (rows 160-163), matching the newly uploaded data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 160: sgt -> sangat ("sangat" already exists as a shared string)
$ws.Range("A160").Value = "sgt"
$ws.Range("B160").Value = "sangat"

# Row 161: skrng -> sekarang
$ws.Range("A161").Value = "skrng"
$ws.Range("B161").Value = "sekarang"

# Row 162: mahall -> mahal
$ws.Range("A162").Value = "mahall"
$ws.Range("B162").Value = "mahal"

# Row 163: kerenn -> keren
# NOTE: "keren" is written before "kerenn" so the shared-string table gets
# the new unique strings in the same order as the source workbook.
$ws.Range("B163").Value = "keren"
$ws.Range("A163").Value = "kerenn"

# Match the author's final view state: scrolled so row 162 is at the top,
# with B164 (the next empty row) selected/active.
$ws.Range("B164").Select()
$excel.ActiveWindow.ScrollRow = 162
$excel.ActiveWindow.ScrollColumn = 1
